$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ 'D' = '25.934.00'; 'E' = '  +0.24%  ' }
    3 = @{ 'D' = '1.641.77'; 'E' = '  +0.20%  ' }
    4 = @{ 'E' = '  -0.03%  ' }
    5 = @{ 'D' = '215.20'; 'E' = '  +0.09%  ' }
    6 = @{ 'D' = '0.5061'; 'E' = '  +0.71%  ' }
    7 = @{ 'D' = '1.002'; 'E' = '  +0.07%  ' }
    8 = @{ 'D' = '0.2560'; 'E' = '  -0.46%  ' }
    9 = @{ 'D' = '0.06369'; 'E' = '  -0.05%  ' }
    10 = @{ 'D' = '19.44' }
    11 = @{ 'D' = '0.07764'; 'E' = '  +0.31%  ' }
    12 = @{ 'D' = '1.654.08'; 'E' = '  +0.85%  ' }
    13 = @{ 'D' = '4.277'; 'E' = '  +0.56%  ' }
    14 = @{ 'D' = '0.5439' }
    15 = @{ 'D' = '0.0₅7818'; 'E' = '  -0.80%  ' }
    16 = @{ 'D' = '64.28'; 'E' = '  +0.22%  ' }
    17 = @{ 'D' = '25.987.45'; 'E' = '  +0.36%  ' }
    18 = @{ 'E' = '  +0.04%  ' }
    19 = @{ 'D' = '197.04'; 'E' = '  -2.31%  ' }
    20 = @{ 'D' = '4.429'; 'E' = '  +1.18%  ' }
    21 = @{ 'D' = '9.935'; 'E' = '  +0.73%  ' }
    22 = @{ 'D' = '6.041'; 'E' = '  +1.22%  ' }
    23 = @{ 'D' = '1.005'; 'E' = '  +0.28%  ' }
    24 = @{ 'D' = '1.895'; 'E' = '  +1.74%  ' }
    25 = @{ 'D' = '140.65'; 'E' = '  -0.02%  ' }
    26 = @{ 'D' = '0.1167'; 'E' = '  +3.02%  ' }
    27 = @{ 'D' = '6.877'; 'E' = '  +1.65%  ' }
    28 = @{ 'E' = '  +0.31%  ' }
    29 = @{ 'D' = '1.236'; 'E' = '  -0.45%  ' }
    30 = @{ 'D' = '0.04955'; 'E' = '  -0.28%  ' }
    31 = @{ 'D' = '3.255'; 'E' = '  -0.37%  ' }
    32 = @{ 'D' = '3.179'; 'E' = '  -0.43%  ' }
    33 = @{ 'E' = '  -0.56%  ' }
    34 = @{ 'D' = '2.363'; 'E' = '  +0.06%  ' }
    35 = @{ 'D' = '0.8943'; 'E' = '  +0.45%  ' }
    36 = @{ 'D' = '2.589'; 'E' = '  -1.48%  ' }
    37 = @{ 'D' = '1.134.16'; 'E' = '  -1.08%  ' }
    38 = @{ 'D' = '0.5441'; 'E' = '  -2.62%  ' }
    39 = @{ 'E' = '  -0.40%  ' }
    40 = @{ 'E' = '  +0.23%  ' }
    41 = @{ 'D' = '2.540'; 'E' = '  -0.64%  ' }
    42 = @{ 'D' = '0.8193'; 'E' = '  +1.93%  ' }
    43 = @{ 'D' = '5.580'; 'E' = '  -1.77%  ' }
    44 = @{ 'D' = '0.0₈126'; 'E' = '  +7.20%  ' }
    45 = @{ 'D' = '99.49'; 'E' = '  -0.09%  ' }
    46 = @{ 'D' = '1.777.37'; 'E' = '  +0.04%  ' }
    47 = @{ 'D' = '0.4534'; 'E' = '  +0.14%  ' }
    48 = @{ 'E' = '  -0.64%  ' }
    49 = @{ 'D' = '54.63'; 'E' = '  -0.01%  ' }
    50 = @{ 'D' = '0.05069'; 'E' = '  +0.33%  ' }
    51 = @{ 'D' = '1.006'; 'E' = '  +0.45%  ' }
}

foreach ($row in $updates.Keys) {
    $rowVals = $updates[$row]
    foreach ($col in $rowVals.Keys) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $rowVals[$col]
        $cell.Style = 'Normal'
    }
}
